$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.194072008132935
$ws.Range("B1").Value = 2.353518724441528
$ws.Range("C1").Value = 3.617908000946045
$ws.Range("D1").Value = 3.183145761489868
$ws.Range("E1").Value = 1.136903882026672
